# Scheduled-runner refresh of the per-sheet Leve profit figures
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4395.778
$ws.Range("I15").Value = 4395.778
$ws.Range("K15").Value = 13187.334
$ws.Range("M15").Value = -13018.334
$ws.Range("H62").Value = 17782266
$ws.Range("I62").Value = 24246792
$ws.Range("K62").Value = 24246792
$ws.Range("M62").Value = -24246168
$ws.Range("H65").Value = 17782266
$ws.Range("I65").Value = 24246792
$ws.Range("K65").Value = 121233960
$ws.Range("M65").Value = -121230840
$ws.Range("H106").Value = 3450
$ws.Range("I106").Value = 3300
$ws.Range("K106").Value = 3300
$ws.Range("M106").Value = -2669
$ws.Range("H133").Value = 118989.664
$ws.Range("J133").Value = 118989.664
$ws.Range("L133").Value = 118989.664
$ws.Range("N133").Value = -129109.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40220.176
$ws.Range("I32").Value = 40349.6
$ws.Range("K32").Value = 40349.6
$ws.Range("M32").Value = -40062.6
$ws.Range("H61").Value = 22234712
$ws.Range("I61").Value = 47628756
$ws.Range("J61").Value = 14924.375
$ws.Range("K61").Value = 47628756
$ws.Range("L61").Value = 14924.375
$ws.Range("M61").Value = -47628544
$ws.Range("N61").Value = -15348.375
$ws.Range("H110").Value = 7814144
$ws.Range("I110").Value = 10870582
$ws.Range("K110").Value = 10870582
$ws.Range("M110").Value = -10868537
$ws.Range("H122").Value = 1622.2
$ws.Range("I122").Value = 1470.6666
$ws.Range("J122").Value = 1849.5
$ws.Range("K122").Value = 4411.9998
$ws.Range("L122").Value = 5548.5
$ws.Range("M122").Value = -1961.9998
$ws.Range("N122").Value = -10448.5
$ws.Range("H132").Value = 3394006.2
$ws.Range("I132").Value = 5002677.5
$ws.Range("K132").Value = 15008032.5
$ws.Range("M132").Value = -15005502.5
$ws.Range("H136").Value = 22234712
$ws.Range("I136").Value = 47628756
$ws.Range("J136").Value = 14924.375
$ws.Range("K136").Value = 142886268
$ws.Range("L136").Value = 44773.125
$ws.Range("M136").Value = -142883718
$ws.Range("N136").Value = -49873.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2886.5833
$ws.Range("I20").Value = 2896.6667
$ws.Range("K20").Value = 2896.6667
$ws.Range("M20").Value = -2649.6667
$ws.Range("H95").Value = 11541.333
$ws.Range("J95").Value = 12312
$ws.Range("L95").Value = 12312
$ws.Range("N95").Value = -17804
$ws.Range("H99").Value = 3032.5
$ws.Range("I99").Value = 2925
$ws.Range("K99").Value = 2925
$ws.Range("M99").Value = -1427
$ws.Range("H105").Value = 66684890
$ws.Range("I105").Value = 90932720
$ws.Range("K105").Value = 90932720
$ws.Range("M105").Value = -90930973

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1418.3077
$ws.Range("I16").Value = 1139.6842
$ws.Range("K16").Value = 1139.6842
$ws.Range("M16").Value = -852.6841999999999
$ws.Range("H31").Value = 35719110
$ws.Range("I31").Value = 90911710
$ws.Range("J31").Value = 6250.294
$ws.Range("K31").Value = 90911710
$ws.Range("L31").Value = 6250.294
$ws.Range("M31").Value = -90911415
$ws.Range("N31").Value = -6840.294
$ws.Range("H34").Value = 35719110
$ws.Range("I34").Value = 90911710
$ws.Range("J34").Value = 6250.294
$ws.Range("K34").Value = 90911710
$ws.Range("L34").Value = 6250.294
$ws.Range("M34").Value = -90911508
$ws.Range("N34").Value = -6654.294
$ws.Range("H113").Value = 1418.3077
$ws.Range("I113").Value = 1139.6842
$ws.Range("K113").Value = 1139.6842
$ws.Range("M113").Value = 1030.3158
$ws.Range("H132").Value = 86831.86
$ws.Range("I132").Value = 6222.5
$ws.Range("K132").Value = 18667.5
$ws.Range("M132").Value = -16137.5
$ws.Range("H141").Value = 454999.8
$ws.Range("J141").Value = 549999.75
$ws.Range("L141").Value = 549999.75
$ws.Range("N141").Value = -560359.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("K32").Value = 1500
$ws.Range("M32").Value = -1217
$ws.Range("H37").Value = 92712.8
$ws.Range("J37").Value = 92712.8
$ws.Range("L37").Value = 278138.4
$ws.Range("N37").Value = -278362.4
$ws.Range("H46").Value = 50
$ws.Range("I46").Value = 50
$ws.Range("K46").Value = 150
$ws.Range("M46").Value = -59
$ws.Range("H69").Value = 2824
$ws.Range("I69").Value = 2040
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 6120
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -5309
$ws.Range("N69").Value = -13622
$ws.Range("H72").Value = 2824
$ws.Range("I72").Value = 2040
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 18360
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -14304
$ws.Range("N72").Value = -44112

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3098.4707
$ws.Range("I80").Value = 2427.5715
$ws.Range("K80").Value = 2427.5715
$ws.Range("M80").Value = -1429.5715
$ws.Range("H83").Value = 3098.4707
$ws.Range("I83").Value = 2427.5715
$ws.Range("K83").Value = 12137.8575
$ws.Range("M83").Value = -7145.8575
$ws.Range("H122").Value = 1959.6
$ws.Range("I122").Value = 2137
$ws.Range("K122").Value = 6411
$ws.Range("M122").Value = -3961
$ws.Range("H139").Value = 100001.336
$ws.Range("J139").Value = 100001.336
$ws.Range("L139").Value = 100001.336
$ws.Range("N139").Value = -110281.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 12000
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H95").Value = 66924
$ws.Range("J95").Value = 66924
$ws.Range("L95").Value = 66924
$ws.Range("N95").Value = -72416
$ws.Range("H96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -50492
$ws.Range("H97").Value = 3171.5
$ws.Range("J97").Value = 3171.5
$ws.Range("L97").Value = 3171.5
$ws.Range("N97").Value = -5153.5
$ws.Range("H122").Value = 3527.7693
$ws.Range("I122").Value = 3385.889
$ws.Range("J122").Value = 3847
$ws.Range("K122").Value = 10157.667
$ws.Range("L122").Value = 11541
$ws.Range("M122").Value = -7707.667000000001
$ws.Range("N122").Value = -16441
$ws.Range("H132").Value = 4846.5
$ws.Range("I132").Value = 2773.8
$ws.Range("J132").Value = 5998
$ws.Range("K132").Value = 8321.400000000001
$ws.Range("L132").Value = 17994
$ws.Range("M132").Value = -5791.400000000001
$ws.Range("N132").Value = -23054

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4320.423
$ws.Range("I81").Value = 4013.28
$ws.Range("J81").Value = 11999
$ws.Range("K81").Value = 8026.56
$ws.Range("L81").Value = 23998
$ws.Range("M81").Value = -6965.56
$ws.Range("N81").Value = -26120
$ws.Range("H84").Value = 4320.423
$ws.Range("I84").Value = 4013.28
$ws.Range("J84").Value = 11999
$ws.Range("K84").Value = 40132.8
$ws.Range("L84").Value = 119990
$ws.Range("M84").Value = -34828.8
$ws.Range("N84").Value = -130598
$ws.Range("H100").Value = 1377.3636
$ws.Range("I100").Value = 1291.75
$ws.Range("J100").Value = 1605.6666
$ws.Range("K100").Value = 2583.5
$ws.Range("L100").Value = 3211.3332
$ws.Range("M100").Value = -2042.5
$ws.Range("N100").Value = -4293.3332
$ws.Range("H101").Value = 25346.834
$ws.Range("J101").Value = 25346.834
$ws.Range("L101").Value = 25346.834
$ws.Range("N101").Value = -31836.834
$ws.Range("H136").Value = 3115.0334
$ws.Range("I136").Value = 1610.9166
$ws.Range("J136").Value = 9131.5
$ws.Range("K136").Value = 4832.7498
$ws.Range("L136").Value = 27394.5
$ws.Range("M136").Value = -2282.7498
$ws.Range("N136").Value = -32494.5
